$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value to a cell while forcing it to be stored as TEXT,
# even when the value looks like a plain number (e.g. "563.38").
# Directly setting .Value on such strings causes automatic numeric
# coercion, so instead we stage the text in a scratch cell that has an
# explicit text NumberFormat, copy it, and paste-special just the value
# into the destination. This keeps the destination cell's own style
# untouched (no "s" attribute gets added) while preserving the text type.
$scratch = $ws.Range("ZZ1")
function Set-TextValue([string]$cellAddr, [string]$val) {
    $scratch.NumberFormat = "@"
    $scratch.Value = $val
    $scratch.Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)
}

$ws.Range("D2").Value = '59.163.73'
$ws.Range("E2").Value = '  -2.32%  '
$ws.Range("D3").Value = '2.581.25'
$ws.Range("E3").Value = '  -2.94%  '
$ws.Range("E4").Value = '  +0.02%  '
Set-TextValue "D5" '563.38'
$ws.Range("E5").Value = '  -0.90%  '
Set-TextValue "D6" '142.33'
$ws.Range("E6").Value = '  -2.07%  '
$ws.Range("E7").Value = '  -0.29%  '
Set-TextValue "D8" '0.597'
$ws.Range("E8").Value = '  -1.50%  '
$ws.Range("D9").Value = '2.588.62'
$ws.Range("E9").Value = '  -2.64%  '
Set-TextValue "D10" '6.64'
$ws.Range("E10").Value = '  -2.77%  '
$ws.Range("E11").Value = '  -0.86%  '
$ws.Range("E12").Value = '  +11.42%  '
$ws.Range("E13").Value = '  +1.72%  '
$ws.Range("D14").Value = '3.037.39'
$ws.Range("E14").Value = '  -2.19%  '
$ws.Range("D15").Value = '59.130.18'
$ws.Range("E15").Value = '  -2.27%  '
Set-TextValue "D16" '22.98'
$ws.Range("E16").Value = '  +5.80%  '
$ws.Range("E17").Value = '  +0.57%  '
$ws.Range("D18").Value = '2.583.96'
$ws.Range("E18").Value = '  -2.36%  '
$ws.Range("E19").Value = '  -0.26%  '
Set-TextValue "D20" '336.96'
$ws.Range("E20").Value = '  -1.82%  '
$ws.Range("E21").Value = '  -0.16%  '
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("E23").Value = '  -0.02%  '
Set-TextValue "D24" '64.07'
$ws.Range("E24").Value = '  -4.07%  '
Set-TextValue "D25" '0.463'
$ws.Range("E25").Value = '  +5.89%  '
Set-TextValue "D26" '0.996'
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("E27").Value = '  -2.90%  '
$ws.Range("E28").Value = '  -0.11%  '
$ws.Range("D29").Value = '0.0₃0775'
$ws.Range("E29").Value = '  +0.57%  '
$ws.Range("E30").Value = '  -0.04%  '
Set-TextValue "D31" '161.70'
$ws.Range("E31").Value = '  +3.48%  '
$ws.Range("E32").Value = '  -2.37%  '
$ws.Range("E33").Value = '  -0.63%  '
Set-TextValue "D34" '18.94'
$ws.Range("E34").Value = '  -1.14%  '
Set-TextValue "D35" '4.03'
$ws.Range("E35").Value = '  -1.39%  '
Set-TextValue "D36" '1.16'
$ws.Range("E36").Value = '  -0.50%  '
Set-TextValue "D37" '0.871'
$ws.Range("E37").Value = '  -3.70%  '
Set-TextValue "D38" '0.875'
$ws.Range("E38").Value = '  -3.75%  '
$ws.Range("E39").Value = '  -0.15%  '
$ws.Range("E40").Value = '  -1.47%  '
Set-TextValue "D41" '294.71'
$ws.Range("E41").Value = '  -3.15%  '
Set-TextValue "D42" '3.66'
$ws.Range("E42").Value = '  -0.02%  '
$ws.Range("E43").Value = '  -0.08%  '
Set-TextValue "D44" '132.00'
$ws.Range("E44").Value = '  +5.95%  '
Set-TextValue "D45" '0.0972'
$ws.Range("E45").Value = '  -0.41%  '
Set-TextValue "D46" '0.596'
$ws.Range("E46").Value = '  -1.46%  '
$ws.Range("B47").Value = 'Hedera'
$ws.Range("C47").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D47" '0.0534'
$ws.Range("E47").Value = '  -2.45%  '
$ws.Range("B48").Value = 'WhiteBITCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue "D48" '10.62'
$ws.Range("E48").Value = '  -0.30%  '
Set-TextValue "D49" '19.03'
$ws.Range("E49").Value = '  -1.64%  '
$ws.Range("E50").Value = '  -0.73%  '
Set-TextValue "D51" '18.54'
$ws.Range("E51").Value = '  +0.64%  '

# Remove the scratch cell so it doesn't affect the sheet's used range.
$scratch.Clear()

